$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.699.43"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "2.674.14"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.89"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.81"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  +2.15%  "
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.92"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.30"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "3.151.64"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "65.496.95"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "2.679.07"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.50"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.32"
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.12"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.76"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("E26").Value = "  -4.03%  "
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "537.25"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  -3.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.75"
$ws.Range("E33").Value = "  -5.69%  "
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.41"
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.423"
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.34"
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.95"
$ws.Range("E40").Value = "  -3.79%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.38"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0610"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.97"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.23"
$ws.Range("E47").Value = "  -5.76%  "
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0997"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("E51").Value = "  +0.21%  "
